$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (e.g. 2021-05-02) updated values
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 20.66443398669174

# Row 3 (e.g. 2021-03-18) updated values
$ws.Range("B3").Value = 0.6545652718822623
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.038307959104277
